# feat: add 2022-Q4 data
#
# Original workbook:
#   Sheet 1: "总计"     -> summary table (one row: 2022-Q3)
#   Sheet 2: "2022-Q3"  -> fund holdings detail for 2022-Q3
#
# Target workbook:
#   Sheet 1: "总计"     -> summary table now has TWO rows: 2022-Q4 (new) then 2022-Q3
#   Sheet 2: "2022-Q4"  -> NEW fund holdings detail for 2022-Q4 (replaces in-place)
#   Sheet 3: "2022-Q3"  -> fund holdings detail for 2022-Q3 (moved here, unchanged)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q4 and
#    push the existing 2022-Q3 row down to row 3.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$wsTotal.Rows(2).Insert()
$wsTotal.Range("A2:D2").ClearFormats()

# Re-apply the same style the old row used (copy format from row 3, which
# still carries the original formatting of column A).
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value2 = 0
$wsTotal.Range("B2").Value2 = "2022-Q4"
$wsTotal.Range("C2").Value2 = 3
$wsTotal.Range("D2").Value2 = 0.02

# The row that used to be "row 2" (2022-Q3) becomes row 3; its index column
# moves from 0 to 1.
$wsTotal.Range("A3").Value2 = 1

# ---------------------------------------------------------------------------
# 2. Duplicate the existing "2022-Q3" detail sheet so the original data is
#    preserved on its own tab, then repurpose the original tab for the new
#    2022-Q4 data (this keeps the sheetId of the "2022-Q3" sheet id=2 name
#    reused for 2022-Q4, with a brand-new sheet id=3 holding the old data).
# ---------------------------------------------------------------------------
$wsQ3Old = $wb.Worksheets.Item(2)

# Full duplicate (values, formats, sheet properties) placed right after it.
$wsQ3Old.Copy($null, $wsQ3Old)
$wsQ3New = $wb.Worksheets.Item(3)

# Rename: old tab becomes 2022-Q4, the duplicate keeps the 2022-Q3 name.
$wsQ3Old.Name = "2022-Q4"
$wsQ3New.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 3. Replace the contents of the (renamed) 2022-Q4 sheet with the new fund
#    holdings data, matching the header style already used elsewhere in the
#    workbook (copied from the "总计" sheet's header cell).
# ---------------------------------------------------------------------------
$wsQ4 = $wsQ3Old

# Drop the old 6-row data body; only 3 data rows are needed for 2022-Q4.
$wsQ4.Rows("5:7").Delete()

# Match the header/index-column style used on the "总计" sheet (style "2" in
# the original file) instead of the style the 2022-Q3 sheet used (style "1").
$wsTotal.Range("B1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2:A4").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsQ4.Range($cols[$i] + "1").Value2 = $headers[$i]
}

function Set-TextCell($ws, $addr, $text) {
    # Force literal-text storage (matches the source data, which stores
    # these numeric-looking values as text, not numbers).
    $ws.Range($addr).Value2 = "'" + $text
}

# Row 2: 519097 / 新华中小市值优选混合
$wsQ4.Range("A2").Value2 = 0
Set-TextCell $wsQ4 "B2" "519097"
Set-TextCell $wsQ4 "C2" "新华中小市值优选混合"
Set-TextCell $wsQ4 "D2" "0.66"
Set-TextCell $wsQ4 "E2" "70.51"
Set-TextCell $wsQ4 "F2" "3.28"
Set-TextCell $wsQ4 "G2" "0.0216"
$wsQ4.Range("H2").Value2 = 7

# Row 3: 013599 / 华润元大臻选回报混合C
$wsQ4.Range("A3").Value2 = 1
Set-TextCell $wsQ4 "B3" "013599"
Set-TextCell $wsQ4 "C3" "华润元大臻选回报混合C"
Set-TextCell $wsQ4 "D3" "0.73"
Set-TextCell $wsQ4 "E3" "62.70"
Set-TextCell $wsQ4 "F3" "0.19"
Set-TextCell $wsQ4 "G3" "0.0014"
$wsQ4.Range("H3").Value2 = 10

# Row 4: 013598 / 华润元大臻选回报混合A
$wsQ4.Range("A4").Value2 = 2
Set-TextCell $wsQ4 "B4" "013598"
Set-TextCell $wsQ4 "C4" "华润元大臻选回报混合A"
Set-TextCell $wsQ4 "D4" "0.00"
Set-TextCell $wsQ4 "E4" "62.70"
Set-TextCell $wsQ4 "F4" "0.19"
$wsQ4.Range("G4").Value2 = 0
$wsQ4.Range("H4").Value2 = 10
